# This script applies a row-level reordering of the weekly price data
# (rows 2-37) in the active worksheet. The underlying per-row records
# (Fecha, Calidad, Volumen, Precio minimo/maximo/promedio, Precio $/Kg)
# are the same set of records as before; only their row positions change
# (data re-sorted/re-shuffled), matching the commit "Fruta / hortaliza, semanal".
#
# Columns A, B, C, E-K, Q, R, T are identical on every data row already,
# so they do not need to be touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44382
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 35000
$ws.Range("O2").Value = 35000
$ws.Range("P2").Value = 35000
$ws.Range("S2").Value = 1944
$ws.Range("D3").Value = 44382
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 32000
$ws.Range("O3").Value = 32000
$ws.Range("P3").Value = 32000
$ws.Range("S3").Value = 1778
$ws.Range("D4").Value = 44382
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = 30000
$ws.Range("O4").Value = 30000
$ws.Range("P4").Value = 30000
$ws.Range("S4").Value = 1667
$ws.Range("D5").Value = 44305
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 25
$ws.Range("N5").Value = 30000
$ws.Range("O5").Value = 30000
$ws.Range("P5").Value = 30000
$ws.Range("S5").Value = 1667
$ws.Range("D6").Value = 44403
$ws.Range("L6").Value = "Especial"
$ws.Range("M6").Value = 25
$ws.Range("N6").Value = 33000
$ws.Range("O6").Value = 33000
$ws.Range("P6").Value = 33000
$ws.Range("S6").Value = 1833
$ws.Range("D7").Value = 44403
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 45
$ws.Range("N7").Value = 30000
$ws.Range("O7").Value = 30000
$ws.Range("P7").Value = 30000
$ws.Range("S7").Value = 1667
$ws.Range("D8").Value = 44403
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 15
$ws.Range("N8").Value = 28000
$ws.Range("O8").Value = 28000
$ws.Range("P8").Value = 28000
$ws.Range("S8").Value = 1556
$ws.Range("D9").Value = 44326
$ws.Range("L9").Value = "Especial"
$ws.Range("M9").Value = 16
$ws.Range("N9").Value = 35000
$ws.Range("O9").Value = 35000
$ws.Range("P9").Value = 35000
$ws.Range("S9").Value = 1944
$ws.Range("D10").Value = 44326
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 25
$ws.Range("N10").Value = 30000
$ws.Range("O10").Value = 30000
$ws.Range("P10").Value = 30000
$ws.Range("S10").Value = 1667
$ws.Range("D11").Value = 44326
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 20
$ws.Range("N11").Value = 28000
$ws.Range("O11").Value = 28000
$ws.Range("P11").Value = 28000
$ws.Range("S11").Value = 1556
$ws.Range("D12").Value = 44424
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 30
$ws.Range("N12").Value = 32000
$ws.Range("O12").Value = 32000
$ws.Range("P12").Value = 32000
$ws.Range("S12").Value = 1778
$ws.Range("D13").Value = 44396
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 35
$ws.Range("N13").Value = 37000
$ws.Range("O13").Value = 37000
$ws.Range("P13").Value = 37000
$ws.Range("S13").Value = 2056
$ws.Range("D14").Value = 44396
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 15
$ws.Range("N14").Value = 34000
$ws.Range("O14").Value = 34000
$ws.Range("P14").Value = 34000
$ws.Range("S14").Value = 1889
$ws.Range("D15").Value = 44431
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 30
$ws.Range("N15").Value = 32000
$ws.Range("O15").Value = 32000
$ws.Range("P15").Value = 32000
$ws.Range("S15").Value = 1778
$ws.Range("D16").Value = 44435
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = 32000
$ws.Range("O16").Value = 32000
$ws.Range("P16").Value = 32000
$ws.Range("S16").Value = 1778
$ws.Range("D17").Value = 44445
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 30
$ws.Range("N17").Value = 32000
$ws.Range("O17").Value = 32000
$ws.Range("P17").Value = 32000
$ws.Range("S17").Value = 1778
$ws.Range("D18").Value = 44319
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 140
$ws.Range("N18").Value = 27000
$ws.Range("O18").Value = 27000
$ws.Range("P18").Value = 27000
$ws.Range("S18").Value = 1500
$ws.Range("D19").Value = 44333
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 30
$ws.Range("N19").Value = 38000
$ws.Range("O19").Value = 38000
$ws.Range("P19").Value = 38000
$ws.Range("S19").Value = 2111
$ws.Range("D20").Value = 44333
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 25
$ws.Range("N20").Value = 35000
$ws.Range("O20").Value = 35000
$ws.Range("P20").Value = 35000
$ws.Range("S20").Value = 1944
$ws.Range("D21").Value = 44284
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 120
$ws.Range("N21").Value = 23000
$ws.Range("O21").Value = 23000
$ws.Range("P21").Value = 23000
$ws.Range("S21").Value = 1278
$ws.Range("D22").Value = 44417
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 15
$ws.Range("N22").Value = 28000
$ws.Range("O22").Value = 28000
$ws.Range("P22").Value = 28000
$ws.Range("S22").Value = 1556
$ws.Range("D23").Value = 44340
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 35
$ws.Range("N23").Value = 37000
$ws.Range("O23").Value = 37000
$ws.Range("P23").Value = 37000
$ws.Range("S23").Value = 2056
$ws.Range("D24").Value = 44340
$ws.Range("L24").Value = "Segunda"
$ws.Range("M24").Value = 20
$ws.Range("N24").Value = 35000
$ws.Range("O24").Value = 35000
$ws.Range("P24").Value = 35000
$ws.Range("S24").Value = 1944
$ws.Range("D25").Value = 44452
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 20
$ws.Range("N25").Value = 36000
$ws.Range("O25").Value = 36000
$ws.Range("P25").Value = 36000
$ws.Range("S25").Value = 2000
$ws.Range("D26").Value = 44354
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 40
$ws.Range("N26").Value = 38000
$ws.Range("O26").Value = 38000
$ws.Range("P26").Value = 38000
$ws.Range("S26").Value = 2111
$ws.Range("D27").Value = 44389
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 35
$ws.Range("N27").Value = 29000
$ws.Range("O27").Value = 29000
$ws.Range("P27").Value = 29000
$ws.Range("S27").Value = 1611
$ws.Range("D28").Value = 44389
$ws.Range("L28").Value = "Segunda"
$ws.Range("M28").Value = 20
$ws.Range("N28").Value = 27000
$ws.Range("O28").Value = 27000
$ws.Range("P28").Value = 27000
$ws.Range("S28").Value = 1500
$ws.Range("D29").Value = 44270
$ws.Range("L29").Value = "Especial"
$ws.Range("M29").Value = 70
$ws.Range("N29").Value = 38000
$ws.Range("O29").Value = 38000
$ws.Range("P29").Value = 38000
$ws.Range("S29").Value = 2111
$ws.Range("D30").Value = 44312
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 160
$ws.Range("N30").Value = 26000
$ws.Range("O30").Value = 26000
$ws.Range("P30").Value = 26000
$ws.Range("S30").Value = 1444
$ws.Range("D31").Value = 44277
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 100
$ws.Range("N31").Value = 30000
$ws.Range("O31").Value = 30000
$ws.Range("P31").Value = 30000
$ws.Range("S31").Value = 1667
$ws.Range("D32").Value = 44277
$ws.Range("L32").Value = "Segunda"
$ws.Range("M32").Value = 60
$ws.Range("N32").Value = 28000
$ws.Range("O32").Value = 28000
$ws.Range("P32").Value = 28000
$ws.Range("S32").Value = 1556
$ws.Range("D33").Value = 44291
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 30
$ws.Range("N33").Value = 25000
$ws.Range("O33").Value = 25000
$ws.Range("P33").Value = 25000
$ws.Range("S33").Value = 1389
$ws.Range("D34").Value = 44438
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 30
$ws.Range("N34").Value = 32000
$ws.Range("O34").Value = 32000
$ws.Range("P34").Value = 32000
$ws.Range("S34").Value = 1778
$ws.Range("D35").Value = 44410
$ws.Range("L35").Value = "Especial"
$ws.Range("M35").Value = 15
$ws.Range("N35").Value = 32000
$ws.Range("O35").Value = 32000
$ws.Range("P35").Value = 32000
$ws.Range("S35").Value = 1778
$ws.Range("D36").Value = 44410
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 25
$ws.Range("N36").Value = 30000
$ws.Range("O36").Value = 30000
$ws.Range("P36").Value = 30000
$ws.Range("S36").Value = 1667
$ws.Range("D37").Value = 44410
$ws.Range("L37").Value = "Segunda"
$ws.Range("M37").Value = 10
$ws.Range("N37").Value = 28000
$ws.Range("O37").Value = 28000
$ws.Range("P37").Value = 28000
$ws.Range("S37").Value = 1556
